# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Wed Jul 10 05:16:58 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.179.15"
$ws.Range("E2").Value = "  +3.40%  "

$ws.Range("D3").Value = "3.109.99"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'523.15"
$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("D6").Value = "'145.06"
$ws.Range("E6").Value = "  +2.93%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +1.66%  "

$ws.Range("D9").Value = "'7.40"
$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("E10").Value = "  +1.74%  "

$ws.Range("E11").Value = "  +3.98%  "

$ws.Range("D12").Value = "3.650.63"
$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("E13").Value = "  +1.71%  "

$ws.Range("D14").Value = "'27.26"
$ws.Range("E14").Value = "  +7.36%  "

$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").Value = "59.162.56"
$ws.Range("E16").Value = "  +3.31%  "

$ws.Range("D17").Value = "3.116.02"
$ws.Range("E17").Value = "  +1.68%  "

$ws.Range("D18").Value = "'6.23"

$ws.Range("D19").Value = "'13.13"
$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("D20").Value = "'8.30"
$ws.Range("E20").Value = "  +2.33%  "

$ws.Range("D21").Value = "'345.08"
$ws.Range("E21").Value = "  +2.49%  "

$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'0.512"
$ws.Range("E23").Value = "  +2.78%  "

$ws.Range("D24").Value = "'66.02"
$ws.Range("E24").Value = "  +1.04%  "

$ws.Range("E25").Value = "  +1.59%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").Value = "0.0₃0934"
$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("D28").Value = "'6.78"
$ws.Range("E28").Value = "  +5.85%  "

$ws.Range("D29").Value = "'7.33"
$ws.Range("E29").Value = "  +4.41%  "

$ws.Range("D30").Value = "'1.86"
$ws.Range("E30").Value = "  +2.84%  "

$ws.Range("E31").Value = "  +4.88%  "

$ws.Range("D32").Value = "'21.14"
$ws.Range("E32").Value = "  +2.11%  "

$ws.Range("D33").Value = "'155.69"
$ws.Range("E33").Value = "  +1.00%  "

$ws.Range("E34").Value = "  +3.63%  "

$ws.Range("D35").Value = "'6.22"
$ws.Range("E35").Value = "  +6.55%  "

$ws.Range("D36").Value = "'27.45"
$ws.Range("E36").Value = "  +5.36%  "

$ws.Range("D37").Value = "'1.32"
$ws.Range("E37").Value = "  +7.11%  "

$ws.Range("D38").Value = "'0.0689"
$ws.Range("E38").Value = "  +3.02%  "

$ws.Range("D39").Value = "'3.98"
$ws.Range("E39").Value = "  +3.83%  "

$ws.Range("D40").Value = "3.158.17"
$ws.Range("E40").Value = "  +1.79%  "

$ws.Range("D41").Value = "'36.95"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.668"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("E44").Value = "  +6.21%  "

$ws.Range("D45").Value = "2.285.50"
$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("E46").Value = "  +3.29%  "

$ws.Range("D47").Value = "'21.17"
$ws.Range("E47").Value = "  +6.01%  "

$ws.Range("D48").Value = "'0.971"
$ws.Range("E48").Value = "  +2.91%  "

$ws.Range("D49").Value = "'6.03"
$ws.Range("E49").Value = "  +3.81%  "

$ws.Range("D50").Value = "'0.758"
$ws.Range("E50").Value = "  +11.19%  "

$ws.Range("D51").Value = "'264.31"
$ws.Range("E51").Value = "  +11.72%  "
